# W5 Salaries and Tasks - fill in this week's team data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block ---------------------------------------------------
$ws.Range("B3").Value = 43754          # Date -> 10/16/2019 (keeps the cell's existing date style/format)
$ws.Range("B5").Value = 5                    # Total Number of Team Members

# --- Team member names & salaries -----------------------------------
$ws.Range("A8").Value = "Kunaal Sikka"
$ws.Range("A9").Value = "Mina Huh"
$ws.Range("A10").Value = "Vu Nguyen"
$ws.Range("A11").Value = "Nicolas Carmody"
$ws.Range("A12").Value = "Jonas Bokstaller"

$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 100

$ws.Range("B4").Value = "MSR Voice Input"   # Team Name

# --- Totals (formulas, unchanged expressions, now resolving live) ---
$ws.Range("B14").Formula = "=SUM(B8:B12)"
$ws.Range("B15").Formula = "=B5*100-B14"

# --- Tasks completed / tasks to do next ------------------------------
$ws.Range("A19").Value = "Finish Presentation"
$ws.Range("A20").Value = "Choose who presents which part"
$ws.Range("A21").Value = "Practice presentation slides"
$ws.Range("B19").Value = "Process presentation feedback"

# --- Column widths (best-fit-ish, matches the refreshed layout) ------
$ws.Columns.Item(1).ColumnWidth = 37.69921875
$ws.Columns.Item(2).ColumnWidth = 34.69921875

# --- View state: scroll/zoom/selection -------------------------------
$ws.Range("B20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 132

# --- Cosmetic: built-in "Normal" cell style is relabeled by a
#     German-locale Excel install that re-saved the file ("Standard").
#     Best-effort; a no-op on hosts that don't expose style renaming.
$wb.Styles.Item("Normal").Name = "Standard"
